$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Cells.Item(10, 2).Value = "SingleUseId10"
$ws.Cells.Item(10, 3).Value = "Label"
$ws.Cells.Item(10, 4).Value = "Center"
$ws.Cells.Item(10, 5).Value = "Send TCP"
$ws.Cells.Item(10, 6).Value = "LTR"

$ws.Range("B10:F10").Style = "Normal"
